# Daily attendance processing - 2025-10-19 13:04:32
#
# Normalizes the "Recorded By" column (G) so that email-address entries
# are listed before plain-name entries (e.g. "System") within each
# comma-separated cell, preserving the relative order within each group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 160   # sheet dimension is A1:S160

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -eq "") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $emails = @()
    $others = @()
    foreach ($p in $trimmed) {
        if ($p.IndexOf("@") -ge 0) {
            $emails += $p
        } else {
            $others += $p
        }
    }

    $ordered = $emails + $others
    $newVal = [string]::Join(", ", $ordered)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
